$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new product rows (63-65) below the existing data.
# Shared-string insertion order matches the authored file: the three
# distinct B-column category labels first, then the repeated C-column
# group label.
$ws.Range("B64").Value = "FAST FOOD VE PAKET SERVİS HİZMETLERİ"
$ws.Range("B65").Value = "KANTİNLER"
$ws.Range("B63").Value = "LOKANTALAR, KAFELER, BARLAR, GECE KULÜPLERİ VB"

$ws.Range("C63").Value = "Restoran ve Kafe"
$ws.Range("C64").Value = "Restoran ve Kafe"
$ws.Range("C65").Value = "Restoran ve Kafe"

$ws.Range("A63").Value = 11111
$ws.Range("A64").Value = 11112
$ws.Range("A65").Value = 11120

# Reflect the scrolled/selected view state from the authored workbook.
[void]$ws.Range("F72").Select()
